$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "-Relationships is not significant:",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "-not significant",
    2
)
